# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Row -> new F value for "展览" sheet
$sheet1Updates = @{
    3  = 1025
    4  = 13470
    9  = 132
    12 = 40
    14 = 13473
    16 = 596
    19 = 8017
    22 = 145
    25 = 7
    26 = 20
    32 = 172
    33 = 376
}

foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# Row -> new F value for "全部类型" sheet (same data, rows 32/33 shifted to 34/35)
$sheet4Updates = @{
    3  = 1025
    4  = 13470
    9  = 132
    12 = 40
    14 = 13473
    16 = 596
    19 = 8017
    22 = 145
    25 = 7
    26 = 20
    34 = 172
    35 = 376
}

foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}

$wb.Save()
